$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.138.83"
$ws.Range("E2").Value = "  -1.68%  "

$ws.Range("D3").Value = "2.173.29"
$ws.Range("E3").Value = "  -2.48%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.20"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -7.87%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.11"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "36.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -12.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0933"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.64%  "

$ws.Range("E13").Value = "  -1.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.03%  "

$ws.Range("D15").Value = "2.498.59"
$ws.Range("E15").Value = "  -2.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.31"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").Value = "2.188.29"
$ws.Range("E18").Value = "  -1.77%  "

$ws.Range("D19").Value = "41.059.76"
$ws.Range("E19").Value = "  -1.81%  "

$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").Value = "  -1.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.68"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.07"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.40"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.10%  "

$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.80"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.41"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()

$ws.Range("E28").Value = "  -4.68%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -7.18%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.22"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.53%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.72%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.94%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0760"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.49%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.121"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.33%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.04%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.96"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.55"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.14%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0305"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.77%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.21"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.30%  "

$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.27"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.52%  "

$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.48"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.54%  "

$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.94"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.55%  "

$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.36"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.15%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.50"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.12%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0997"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.81%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.189"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.89%  "

$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.14"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.96%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E50").Value = "  -4.66%  "

$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.73"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.07%  "

Write-Output "Applied cryptos update"